$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Open Action Items")

# Shared replacement strings
$depOld = "Dependent on AI/ML Implementation milestone completion"
$depNew = "Dependent on Cloud Infrastructure Migration milestone completion"
$noteOld = "Critical action for Artificial Intelligence and Machine Learning success"
$noteNew = "Critical action for Information Technology success"

# Row 8
$ws.Range("E8").Value = "Chief Technology Officer"
$ws.Range("I8").Value = $depNew
$ws.Range("J8").Value = $noteNew

# Row 9
$ws.Range("E9").Value = "IT Managers"
$ws.Range("I9").Value = $depNew
$ws.Range("J9").Value = $noteNew

# Row 10
$ws.Range("E10").Value = "DevOps Engineers"
$ws.Range("I10").Value = $depNew
$ws.Range("J10").Value = $noteNew

# Row 11
$ws.Range("E11").Value = "System Administrators"
$ws.Range("I11").Value = $depNew
$ws.Range("J11").Value = $noteNew

# Row 12 (Owner unchanged)
$ws.Range("I12").Value = $depNew
$ws.Range("J12").Value = $noteNew

# Row 13 (Owner unchanged)
$ws.Range("I13").Value = $depNew
$ws.Range("J13").Value = $noteNew

# Row 14
$ws.Range("E14").Value = "Chief Technology Officer"
$ws.Range("I14").Value = $depNew
$ws.Range("J14").Value = $noteNew

# Row 15
$ws.Range("E15").Value = "IT Managers"
$ws.Range("I15").Value = $depNew
$ws.Range("J15").Value = $noteNew

# Row 16
$ws.Range("E16").Value = "DevOps Engineers"
$ws.Range("I16").Value = $depNew
$ws.Range("J16").Value = $noteNew

# Row 17
$ws.Range("E17").Value = "System Administrators"
$ws.Range("I17").Value = $depNew
$ws.Range("J17").Value = $noteNew

$wb.Save()
